$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26. This pushes the existing rows 26..160
# down to 27..161 (preserving all their values/formatting), matching the
# target diff where every record from the old row 26 onward shifts down
# by one row and a brand new record appears at row 26.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new price-observation record.
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44677
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112017
$ws.Range("G26").Value = "Apio"
$ws.Range("H26").Value = "Americana (o)"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 8000
$ws.Range("N26").Value = "`$/docena de matas"
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 1333
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = "Hortaliza"

# Ensure the date cell uses the same date/time number format as the other
# rows' date column (style index 2 in the original workbook).
$ws.Range("D26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
